$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Folio No" in G1 (new shared string, extends used range to G)
$ws.Range("G1").Value = "Folio No"

# Match the cursor/selection position shown in the diff (G2)
$ws.Range("G2").Select()
